# Refresh the cryptocurrency price/volume snapshot (GitHub Actions cron update).
# For numeric-looking "Price" values we temporarily force a text number format
# before assigning, then clear the format override, so the value is stored as
# text (matching the sheet's existing convention) instead of being silently
# coerced into a floating point number by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.119.40'
$ws.Range('E2').Value = '  -1.86%  '
$ws.Range('D3').Value = '3.346.23'
$ws.Range('E3').Value = '  -3.55%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '559.46'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.48%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '176.53'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.77%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.615'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.65%  '
$ws.Range('D8').Value = '3.337.19'
$ws.Range('E8').Value = '  -3.63%  '
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  -0.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.162'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +5.94%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '54.11'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +1.30%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000271'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.28%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.07'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.17%  '
$ws.Range('D15').Value = '3.885.21'
$ws.Range('E15').Value = '  -3.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '18.32'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.16%  '
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('B18').Value = 'Uniswap'
$ws.Range('C18').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.91'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = '3.351.21'
$ws.Range('E19').Value = '  -3.69%  '
$ws.Range('D20').Value = '64.124.12'
$ws.Range('E20').Value = '  -1.95%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.984'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.35%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '441.09'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +7.67%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.13'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.50'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +9.02%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.69'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '13.30'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +4.47%  '
$ws.Range('B27').Value = 'ImmutableX'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.85'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +1.85%  '
$ws.Range('B28').Value = 'RenderToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.67'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.81'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -0.53%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '29.64'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -0.88%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.58'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  +4.99%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '592.25'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -4.44%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.51'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.108'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '58.66'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.08%  '
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('E37').Value = '  -0.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.50'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.91%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '35.59'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -2.70%  '
$ws.Range('D40').Value = '0.0₃0755'
$ws.Range('E40').Value = '  -3.42%  '
$ws.Range('E41').Value = '  -1.66%  '
$ws.Range('D42').Value = '3.135.64'
$ws.Range('E42').Value = '  -6.64%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -0.02%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.88'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.85%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0407'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.71%  '
$ws.Range('E46').Value = '  -1.79%  '
$ws.Range('E47').Value = '  -0.99%  '
$ws.Range('E48').Value = '  -1.34%  '
$ws.Range('E49').Value = '  -3.81%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.23'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '134.08'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.14%  '
